# Correzione TraceId test #6 e #147
# - Rename sheet "Test Cases" -> "TestCases" (also updates the
#   _xlnm._FilterDatabase defined name automatically)
# - Bump the sheet zoom level 55% -> 70%
# - Fix the TRACEID values for test #6 (row 8) and test #147 (row 88)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (Excel keeps the FilterDatabase defined name in sync)
$ws.Name = "TestCases"

# Update the zoom level of the sheet view
$excel.ActiveWindow.Zoom = 70

# Test #6 (row 8): TRACEID eae8d69bef653d12 -> d6ba896e12aecce9
$ws.Range("H8").Value = "d6ba896e12aecce9"

# Test #147 (row 88): TRACEID 641ec7d64dc680c2 -> ae72ee4beb45bf2f
$ws.Range("H88").Value = "ae72ee4beb45bf2f"
